$d = $word.ActiveDocument

$urls = @(
  "https://www.amazon.com/Slamtec-RPLIDAR-Scanning-Avoidance-Navigation/dp/B07TJW5SXF/ref=sr_1_3?keywords=lidar&qid=1638967125&sr=8-3",
  "https://arstechnica.com/gadgets/2020/06/boston-dynamics-robot-dog-can-be-yours-for-the-low-low-price-of-74500/",
  "https://www.mathworks.com/videos/matlab-and-simulink-robotics-arena-walking-robots-pattern-generation-1546434170253.html?s_tid=srchtitle_Modeling%20and%20Simulation%20of%20Walking%20Robots_5",
  "https://github.com/mathworks-robotics",
  "https://www.mathworks.com/matlabcentral/fileexchange/64227-matlab-and-simulink-robotics-arena-walking-robot?s_tid=srchtitle_Modeling%20and%20Simulation%20of%20Walking%20Robots_4",
  "https://www.mathworks.com/videos/deep-reinforcement-learning-for-walking-robots--1551449152203.html?s_tid=srchtitle_Modeling%20and%20Simulation%20of%20Walking%20Robots_3",
  "https://blogs.mathworks.com/student-lounge/2019/12/20/walking-robot-modeling-and-simulation/?s_tid=srchtitle_Modeling%20and%20Simulation%20of%20Walking%20Robots_2",
  "https://www.mathworks.com/videos/modeling-and-simulation-of-walking-robots-1576560207573.html?s_tid=srchtitle_Modeling%20and%20Simulation%20of%20Walking%20Robots_1",
  "https://www.mathworks.com/matlabcentral/answers/251947-basic-things-for-creating-a-simulated-robot-walk?s_tid=srchtitle_Modeling%2520and%2520Simulation%2520of%2520Walking%2520Robots_10",
  "https://www.mathworks.com/help/deeplearning/ug/train-biped-robot-to-walk-using-reinforcement-learning-agents.html?s_tid=srchtitle_Modeling%2520and%2520Simulation%2520of%2520Walking%2520Robots_9",
  "https://www.mathworks.com/help/reinforcement-learning/ug/train-biped-robot-to-walk-using-reinforcement-learning-agents.html?s_tid=srchtitle_Modeling%2520and%2520Simulation%2520of%2520Walking%2520Robots_8",
  "https://www.mathworks.com/videos/reinforcement-learning-part-4-the-walking-robot-problem-1557482052319.html?s_tid=srchtitle_Modeling%2520and%2520Simulation%2520of%2520Walking%2520Robots_7",
  "https://www.mathworks.com/videos/matlab-and-simulink-robotics-arena-walking-robots-part-3-trajectory-optimization-1506440520726.html?s_tid=srchtitle_Modeling%20and%20Simulation%20of%20Walking%20Robots_6",
  "https://www.mathworks.com/search.html?c%5B%5D=entire_site&q=Modeling%20and%20Simulation%20of%20Walking%20Robots&page=2"
)

foreach ($url in $urls) {
  $end = $d.Content
  $end.Collapse(0)
  $p = $end.Paragraphs.Add()
  $p.Range.Text = $url
}

Write-Host "Done"
